$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
  "G2" = 39.41161066666667
  "H2" = 118.234832
  "I2" = 0.2026411830173254
  "J2" = 0.2026411830173254
  "M2" = 3.626135
  "N2" = 10.878405
  "O2" = 0.4728835835086186
  "P2" = 0.4728835835086186
  "Q2" = 142.9118208447733
  "R2" = 1286.20638760296
  "S2" = 0.09582568879165866
  "T2" = 0.09582568879165865
  "G3" = 39.41161066666667
  "H3" = 118.234832
  "I3" = 0.2026411830173254
  "J3" = 0.2026411830173254
  "O3" = 0.01581792773244636
  "P3" = 0.01581792773244636
  "Q3" = 4.780391904202667
  "R3" = 43.023527137824
  "S3" = 0.00320536358858549
  "T3" = 0.003205363588585489
  "G4" = 39.41161066666667
  "H4" = 118.234832
  "I4" = 0.2026411830173254
  "J4" = 0.2026411830173254
  "O4" = 0.511298488758935
  "P4" = 0.511298488758935
  "Q4" = 154.5213252732605
  "R4" = 1390.691927459344
  "S4" = 0.1036101306370812
  "T4" = 0.1036101306370812
  "I5" = 0.6376490878685519
  "J5" = 0.6376490878685519
  "M5" = 3.626135
  "N5" = 10.878405
  "O5" = 0.4728835835086186
  "P5" = 0.4728835835086186
  "Q5" = 449.6992706537467
  "R5" = 4047.293435883721
  "S5" = 0.3015337856922828
  "T5" = 0.3015337856922829
  "I6" = 0.6376490878685519
  "J6" = 0.6376490878685519
  "O6" = 0.01581792773244636
  "P6" = 0.01581792773244636
  "S6" = 0.01008628719056509
  "T6" = 0.01008628719056509
  "I7" = 0.6376490878685519
  "J7" = 0.6376490878685519
  "O7" = 0.511298488758935
  "P7" = 0.511298488758935
  "S7" = 0.326029014985704
  "T7" = 0.326029014985704
  "I8" = 0.1597097291141227
  "J8" = 0.1597097291141227
  "M8" = 3.626135
  "N8" = 10.878405
  "O8" = 0.4728835835086186
  "P8" = 0.4728835835086186
  "Q8" = 112.6345980341683
  "R8" = 1013.711382307515
  "S8" = 0.0755241090246771
  "T8" = 0.07552410902467711
  "I9" = 0.1597097291141227
  "J9" = 0.1597097291141227
  "O9" = 0.01581792773244636
  "P9" = 0.01581792773244636
  "S9" = 0.002526276953295777
  "T9" = 0.002526276953295777
  "I10" = 0.1597097291141227
  "J10" = 0.1597097291141227
  "O10" = 0.511298488758935
  "P10" = 0.511298488758935
  "S10" = 0.08165934313614982
  "T10" = 0.08165934313614982
}

foreach ($cell in $updates.Keys) {
  $ws.Range($cell).Value = $updates[$cell]
}
